$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite rows 2-4 with the 2010-2012 data (previously rows 12-14)
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 0.037
$ws.Range("C2").Value = 0.105
$ws.Range("D2").Value = 0.3493
$ws.Range("E2").Value = 2.2792

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 0.0375
$ws.Range("C3").Value = 0.1075
$ws.Range("D3").Value = 0.4879
$ws.Range("E3").Value = 2.301

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 0.035769943
$ws.Range("C4").Value = 0.101728418
$ws.Range("D4").Value = 0.484185071
$ws.Range("E4").Value = 2.339326526

# Remove the now-obsolete rows (previously held years 2002-2012 data
# beyond row 4) so the used range shrinks back down to A1:E4.
$ws.Range("A5:E14").EntireRow.Delete()
